$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week of 31 July - 6 August 2016 timesheet block ---
# Monday (row 39): updated activity text, fixing "Progres" -> "Progress" typo
$ws.Range("B39").Value = "Updated sample code to use threads for background processing" + [char]10 + "Updated DSDM methodology " + [char]10 + "Updated Project Diagram" + [char]10 + "Added sample GUI with Progress Table"

# Tuesday (row 40): updated activity text and hours worked (7 -> 10)
$ws.Range("B40").Value = "Updated Project Analysis" + [char]10 + "Updated Project Proposal" + [char]10 + "Updated sample code for GUI progress table"
$ws.Range("D40").Value = 10

# Wednesday (row 41): updated activity text
$ws.Range("B41").Value = "Updated time box items" + [char]10 + "Updated logos for splash screen and main app" + [char]10 + "Updated sample code which allows disabling of resizing of windows form"

# --- New week header (row 47) ---
$ws.Range("A47").Value = "Week Sunday, 7 August - Saturday, 13 August 2016"

# --- Sunday entry of the new week (row 49): append a new activity line ---
$ws.Range("B49").Value = "Updated time sheet and had project meeting" + [char]10 + "Added sample code for encrypting files"
$ws.Rows.Item(49).RowHeight = 57.75

# --- New grand-total row summing every week's subtotal ---
$ws.Range("D59").Formula = "=SUM(D56,D45,D34,D23,D12)"

# --- Restore selection to match the saved view state ---
$ws.Range("B50").Select()
